$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals = @(3, 1, 0.5396133333333334, 1.61884, 0.02587018426425635, 0.02587018426425635, 1, 0.3333333333333333, 0.04871, 0.14613, 0.06500340741874241, 0.06500340741874241, 0.02628456546666667, 0.2365610892, 0.001681650127727395, 0.001681650127727395)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(2, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 0.5396133333333334, 1.61884, 0.02587018426425635, 0.02587018426425635, 1, 0.3333333333333333, 0.059142, 0.177426, 0.07892489266186128, 0.07892489266186128, 0.03191381176, 0.28722430584, 0.002041801516199006, 0.002041801516199006)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(3, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 0.5396133333333334, 1.61884, 0.02587018426425635, 0.02587018426425635, 2, 0.6666666666666666, 0.3835156666666666, 1.150547, 0.511800967600163, 0.511800967600163, 0.2069501672755556, 1.86255150548, 0.01324038533844091, 0.01324038533844091)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(4, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 0.5396133333333334, 1.61884, 0.02587018426425635, 0.02587018426425635, 3, 1, 0.1956336666666667, 0.586901, 0.2610727764146126, 0.2610727764146126, 0.1055665349822222, 0.95009881484, 0.006754000832227029, 0.006754000832227027)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(5, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 0.5396133333333334, 1.61884, 0.02587018426425635, 0.02587018426425635, 2, 0.6666666666666666, 0.06234399999999999, 0.187032, 0.08319795590462073, 0.08319795590462074, 0.03364165365333333, 0.30277488288, 0.002152346449662013, 0.002152346449662014)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(6, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 6.260434, 18.781302, 0.3001382122153186, 0.3001382122153186, 1, 0.3333333333333333, 0.04871, 0.14613, 0.06500340741874241, 0.06500340741874241, 0.30494574014, 2.74451166126, 0.01951000649056533, 0.01951000649056533)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(7, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 6.260434, 18.781302, 0.3001382122153186, 0.3001382122153186, 1, 0.3333333333333333, 0.059142, 0.177426, 0.07892489266186128, 0.07892489266186128, 0.370254587628, 3.332291288652, 0.02368837618281696, 0.02368837618281696)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(8, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 6.260434, 18.781302, 0.3001382122153186, 0.3001382122153186, 2, 0.6666666666666666, 0.3835156666666666, 1.150547, 0.511800967600163, 0.511800967600163, 2.400974519132666, 21.608770672194, 0.1536110274255831, 0.1536110274255831)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(9, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 6.260434, 18.781302, 0.3001382122153186, 0.3001382122153186, 3, 1, 0.1956336666666667, 0.586901, 0.2610727764146126, 0.2610727764146126, 1.224751658344667, 11.022764925102, 0.07835791637117143, 0.07835791637117141)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(10, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 6.260434, 18.781302, 0.3001382122153186, 0.3001382122153186, 2, 0.6666666666666666, 0.06234399999999999, 0.187032, 0.08319795590462073, 0.08319795590462074, 0.3903004972959999, 3.512704475664, 0.02497088574518178, 0.02497088574518178)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(11, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 1.182122, 3.546366, 0.0566733845769154, 0.0566733845769154, 1, 0.3333333333333333, 0.04871, 0.14613, 0.06500340741874241, 0.06500340741874241, 0.05758116262, 0.51823046358, 0.003683963107452305, 0.003683963107452305)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(12, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 1.182122, 3.546366, 0.0566733845769154, 0.0566733845769154, 1, 0.3333333333333333, 0.059142, 0.177426, 0.07892489266186128, 0.07892489266186128, 0.06991305932399999, 0.629217533916, 0.004472940794517434, 0.004472940794517434)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(13, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 1.182122, 3.546366, 0.0566733845769154, 0.0566733845769154, 2, 0.6666666666666666, 0.3835156666666666, 1.150547, 0.511800967600163, 0.511800967600163, 0.4533623069113333, 4.080260762202, 0.02900549306364146, 0.02900549306364146)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(14, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 1.182122, 3.546366, 0.0566733845769154, 0.0566733845769154, 3, 1, 0.1956336666666667, 0.586901, 0.2610727764146126, 0.2610727764146126, 0.2312628613073333, 2.081365751766, 0.01479587786030839, 0.01479587786030839)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(15, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 1.182122, 3.546366, 0.0566733845769154, 0.0566733845769154, 2, 0.6666666666666666, 0.06234399999999999, 0.187032, 0.08319795590462073, 0.08319795590462074, 0.07369821396799998, 0.6632839257119999, 0.004715109750995821, 0.004715109750995821)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(16, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 12.87633433333333, 38.629003, 0.6173182189435098, 0.6173182189435098, 1, 0.3333333333333333, 0.04871, 0.14613, 0.06500340741874241, 0.06500340741874241, 0.6272062453766667, 5.64485620839, 0.0401277876929974, 0.0401277876929974)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(17, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 12.87633433333333, 38.629003, 0.6173182189435098, 0.6173182189435098, 1, 0.3333333333333333, 0.059142, 0.177426, 0.07892489266186128, 0.07892489266186128, 0.761532165142, 6.853789486277999, 0.0487217741683279, 0.0487217741683279)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(18, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 12.87633433333333, 38.629003, 0.6173182189435098, 0.6173182189435098, 2, 0.6666666666666666, 0.3835156666666666, 1.150547, 0.511800967600163, 0.511800967600163, 4.938275946071221, 44.44448351464099, 0.3159440617724976, 0.3159440617724976)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(19, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 12.87633433333333, 38.629003, 0.6173182189435098, 0.6173182189435098, 3, 1, 0.1956336666666667, 0.586901, 0.2610727764146126, 0.2610727764146126, 2.519044498855889, 22.671400489703, 0.1611649813509058, 0.1611649813509058)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(20, 5 + $i).Value2 = $vals[$i] }
$vals = @(3, 1, 12.87633433333333, 38.629003, 0.6173182189435098, 0.6173182189435098, 2, 0.6666666666666666, 0.06234399999999999, 0.187032, 0.08319795590462073, 0.08319795590462074, 0.8027621876773331, 7.224859689095998, 0.05135961395878114, 0.05135961395878114)
for ($i = 0; $i -lt $vals.Count; $i++) { $ws.Cells.Item(21, 5 + $i).Value2 = $vals[$i] }
Write-Output "applied update"
